# Updated cryptos list - applies latest price/volume snapshot to the sheet,
# including a swap of the Chainlink/BitcoinCash rows (18 and 19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.947.14"
$ws.Range("E2").Value = "  -0.06%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.671.70"
$ws.Range("E3").Value = "  +1.13%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'214.89"
$ws.Range("E5").Value = "  +0.00%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.58%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.05%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +0.27%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +0.56%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'20.19"
$ws.Range("E10").Value = "  +0.06%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.34%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.906.94"
$ws.Range("E12").Value = "  +1.11%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.686.57"
$ws.Range("E13").Value = "  +2.05%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.08%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +1.00%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "'65.54"
$ws.Range("E16").Value = "  +0.50%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.948.30"
$ws.Range("E17").Value = "  -0.06%  "

# Row 18 / Row 19 - Chainlink and BitcoinCash swap places (Chainlink drops to
# row 19, BitcoinCash moves up to row 18), with refreshed price/volume data.
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'234.18"
$ws.Range("E18").Value = "  -0.98%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'8.03"
$ws.Range("E19").Value = "  +3.49%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0732"
$ws.Range("E20").Value = "  -0.12%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.01%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +0.19%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  -1.44%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -2.00%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "'7.13"
$ws.Range("E26").Value = "  +0.15%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'15.95"
$ws.Range("E27").Value = "  +0.76%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -1.32%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.07%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.08%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.42%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.460.54"
$ws.Range("E33").Value = "  -5.86%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +1.77%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +1.76%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  +0.13%  "

# Row 37 - ImmutableX
$ws.Range("D37").Value = "'0.579"
$ws.Range("E37").Value = "  -0.52%  "

# Row 38 - ARBITRUM
$ws.Range("D38").Value = "'0.899"
$ws.Range("E38").Value = "  -0.18%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +1.00%  "

# Row 40 - WEMIXToken
$ws.Range("E40").Value = "  +12.36%  "

# Row 41 - FraxShare
$ws.Range("D41").Value = "'5.79"
$ws.Range("E41").Value = "  -3.30%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.08%  "

# Row 43 - MXToken
$ws.Range("E43").Value = "  +2.80%  "

# Row 44 - Aave
$ws.Range("D44").Value = "'66.53"
$ws.Range("E44").Value = "  +0.15%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.812.77"
$ws.Range("E45").Value = "  +1.04%  "

# Row 46 - TrustWalletToken
$ws.Range("E46").Value = "  +0.99%  "

# Row 47 - Quant
$ws.Range("D47").Value = "'90.70"
$ws.Range("E47").Value = "  +0.98%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  +0.70%  "

# Row 49 - Algorand
$ws.Range("E49").Value = "  +2.68%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  +0.50%  "

# Row 51 - EnergySwap
$ws.Range("E51").Value = "  +0.45%  "
